# NIT-9014706146.xlsx — "Actualiza base de datos EC y agrega parte 1 de
# nuevos estado de cuenta"
#
# The account-statement worksheet is refreshed with a new batch of data:
#   - The stale "mora" total and worker/period counters are updated.
#   - The five now-outdated worker rows (NATALIA BERRIO DUARTE / ELIANNY
#     CLARET SALCEDO RODRIGUEZ / VICTOR MANUEL VILLALOBO) are removed,
#     leaving just the first worker row as "parte 1" of the new statement.
#   - The trailing signature block shifts up to follow the shorter table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the five obsolete worker rows (old rows 17-21); everything below
# (the signature/legal-representative block) shifts up to close the gap.
$ws.Rows("17:21").Delete()

# Refresh the summary figures for the new batch.
$ws.Range("E11").Value = 7592   # VALOR MORA
$ws.Range("C13").Value = 1      # Cant. Trabajadores
$ws.Range("F13").Value = 1      # Cant. Periodos

# Column D ("Nombre Trabajador") no longer needs to fit the longest of the
# removed names, so re-run best-fit sizing now that only the shorter
# "PAOLA ANDREA TORRES ROJAS" remains.
$ws.Columns("D:D").AutoFit()
